$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update SUB_TOTAL for row 8 (C8): 0 -> 0.02 (materials now loading in correctly)
$ws.Range("C8").Value = 0.02

# Add a new note in D8 about needing a README.txt file of controls
$ws.Range("D8").Value = "Need a README.txt file of controls!"

# Update the active selection to E8 (matches author's final cursor position)
$ws.Range("E8").Select()

# Reflect the author's resized/repositioned window, if supported by this host
$win = $excel.ActiveWindow
$win.Left = 3795
$win.Top = 375
$win.Width = 25155
$win.Height = 18090
